$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary rows 14-17: label in column A, bold 12pt vertically-centered
# aggregate formula in column B.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Copy B14's formatting onto B15:B17 (format-painter style) so they share the
# same cell style without generating redundant/orphaned style records.
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Row 12: bold 11pt average of column J.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

$ws.Application.CutCopyMode = $false
$ws.Range("J12").Select() | Out-Null
